# 1. remove outdated interfaces
#    - drop the duplicate/unused cell style (cellXfs[1]) by resetting the
#      number format on the data cells that used it back to "General" so
#      they fall back onto the shared base style.
# 2. weighted configuration selection
#    - add a "job_id" header/identifier column (A) ahead of the existing
#      instance-type columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert the new leading "job_id" label in A1.
$ws.Range("A1").Value = "job_id"

# Normalize the formatting of the data block so it no longer references the
# now-redundant second cell style (merges style 1 back into style 0).
$ws.Range("B2:F6").NumberFormat = "General"

# Match the recorded selection state.
$ws.Range("B8").Select() | Out-Null
